$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file moved one minute later
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 07:22:21"

# zh-cn sheet: handoff / handback datetimes for the first file updated
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 07:22:16"
$wsZhCn.Range("K2").Value = "2016-09-01 07:22:33"

# de-de sheet: handoff / handback datetimes for the first file updated
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 07:22:21"
$wsDeDe.Range("K2").Value = "2016-09-01 07:22:40"
